# Daily attendance processing - normalize "Recorded By" ordering.
# Swap "System, <email>" -> "<email>, System" throughout the sheet
# (column G / "Recorded By"), matching the latest attendance sync output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
